$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(17, 8).Value = 1855.909
$ws.Cells.Item(17, 10).Value = 1855.909
$ws.Cells.Item(17, 12).Value = 5567.727000000001
$ws.Cells.Item(17, 14).Value = -5903.727000000001
$ws.Cells.Item(58, 8).Value = 1895.8
$ws.Cells.Item(58, 9).Value = 1313.3334
$ws.Cells.Item(58, 10).Value = 2769.5
$ws.Cells.Item(58, 11).Value = 3940.0002
$ws.Cells.Item(58, 12).Value = 8308.5
$ws.Cells.Item(58, 13).Value = -3790.0002
$ws.Cells.Item(58, 14).Value = -8608.5
$ws.Cells.Item(86, 8).Value = 9115.333000000001
$ws.Cells.Item(86, 9).Value = 8232.333000000001
$ws.Cells.Item(86, 10).Value = 9998.333000000001
$ws.Cells.Item(86, 11).Value = 8232.333000000001
$ws.Cells.Item(86, 12).Value = 9998.333000000001
$ws.Cells.Item(86, 13).Value = -7109.333000000001
$ws.Cells.Item(86, 14).Value = -12244.333
$ws.Cells.Item(89, 8).Value = 9115.333000000001
$ws.Cells.Item(89, 9).Value = 8232.333000000001
$ws.Cells.Item(89, 10).Value = 9998.333000000001
$ws.Cells.Item(89, 11).Value = 41161.665
$ws.Cells.Item(89, 12).Value = 49991.665
$ws.Cells.Item(89, 13).Value = -35545.665
$ws.Cells.Item(89, 14).Value = -61223.665

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(5, 8).Value = 72.57143000000001
$ws.Cells.Item(5, 9).Value = 61.6
$ws.Cells.Item(5, 11).Value = 61.6
$ws.Cells.Item(5, 13).Value = 50.4
$ws.Cells.Item(61, 8).Value = 6849.1
$ws.Cells.Item(61, 9).Value = 5415.1665
$ws.Cells.Item(61, 10).Value = 9000
$ws.Cells.Item(61, 11).Value = 5415.1665
$ws.Cells.Item(61, 12).Value = 9000
$ws.Cells.Item(61, 13).Value = -5203.1665
$ws.Cells.Item(61, 14).Value = -9424
$ws.Cells.Item(132, 8).Value = 170.25
$ws.Cells.Item(132, 9).Value = 170.25
$ws.Cells.Item(132, 11).Value = 510.75
$ws.Cells.Item(132, 13).Value = 2019.25
$ws.Cells.Item(136, 8).Value = 6849.1
$ws.Cells.Item(136, 9).Value = 5415.1665
$ws.Cells.Item(136, 10).Value = 9000
$ws.Cells.Item(136, 11).Value = 16245.4995
$ws.Cells.Item(136, 12).Value = 27000
$ws.Cells.Item(136, 13).Value = -13695.4995
$ws.Cells.Item(136, 14).Value = -32100

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(4, 8).Value = 72.57143000000001
$ws.Cells.Item(4, 9).Value = 61.6
$ws.Cells.Item(4, 11).Value = 61.6
$ws.Cells.Item(4, 13).Value = 53.4
$ws.Cells.Item(94, 8).Value = 290.85715
$ws.Cells.Item(94, 9).Value = 287.4
$ws.Cells.Item(94, 11).Value = 287.4
$ws.Cells.Item(94, 13).Value = 163.6
$ws.Cells.Item(99, 8).Value = 3100.4443
$ws.Cells.Item(99, 9).Value = 2682
$ws.Cells.Item(99, 10).Value = 3623.5
$ws.Cells.Item(99, 11).Value = 2682
$ws.Cells.Item(99, 12).Value = 3623.5
$ws.Cells.Item(99, 13).Value = -1184
$ws.Cells.Item(99, 14).Value = -6619.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(132, 8).Value = 1954.25
$ws.Cells.Item(132, 9).Value = 1954.25
$ws.Cells.Item(132, 11).Value = 5862.75
$ws.Cells.Item(132, 13).Value = -3332.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 39.545456
$ws.Cells.Item(2, 9).Value = 21
$ws.Cells.Item(2, 10).Value = 50.142857
$ws.Cells.Item(2, 11).Value = 126
$ws.Cells.Item(2, 12).Value = 300.857142
$ws.Cells.Item(2, 13).Value = -13
$ws.Cells.Item(2, 14).Value = -526.8571420000001
$ws.Cells.Item(4, 8).Value = 264092.5
$ws.Cells.Item(4, 10).Value = 36829.668
$ws.Cells.Item(4, 12).Value = 110489.004
$ws.Cells.Item(4, 14).Value = -110713.004
$ws.Cells.Item(22, 8).Value = 293.57144
$ws.Cells.Item(22, 10).Value = 295.83334
$ws.Cells.Item(22, 12).Value = 887.5000200000001
$ws.Cells.Item(22, 14).Value = -1225.50002
$ws.Cells.Item(27, 8).Value = 293.57144
$ws.Cells.Item(27, 10).Value = 295.83334
$ws.Cells.Item(27, 12).Value = 887.5000200000001
$ws.Cells.Item(27, 14).Value = -1091.50002
$ws.Cells.Item(34, 8).Value = 970.9286
$ws.Cells.Item(34, 9).Value = 206.33333
$ws.Cells.Item(34, 10).Value = 1179.4546
$ws.Cells.Item(34, 11).Value = 618.99999
$ws.Cells.Item(34, 12).Value = 3538.3638
$ws.Cells.Item(34, 13).Value = -534.99999
$ws.Cells.Item(34, 14).Value = -3706.3638
$ws.Cells.Item(118, 8).Value = 888.5
$ws.Cells.Item(118, 9).Value = 888.5
$ws.Cells.Item(118, 11).Value = 2665.5
$ws.Cells.Item(118, 13).Value = -1422.5
$ws.Cells.Item(121, 8).Value = 999
$ws.Cells.Item(121, 9).Value = 999
$ws.Cells.Item(121, 10).Value = 0
$ws.Cells.Item(121, 11).Value = 2997
$ws.Cells.Item(121, 12).Value = 0
$ws.Cells.Item(121, 13).Value = -1687
$ws.Cells.Item(121, 14).ClearContents()
$ws.Cells.Item(122, 8).Value = 897.5
$ws.Cells.Item(122, 10).Value = 972.75
$ws.Cells.Item(122, 12).Value = 8754.75
$ws.Cells.Item(122, 14).Value = -13654.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 7166.6665
$ws.Cells.Item(70, 9).Value = 5750
$ws.Cells.Item(70, 10).Value = 10000
$ws.Cells.Item(70, 11).Value = 5750
$ws.Cells.Item(70, 12).Value = 10000
$ws.Cells.Item(70, 13).Value = -5480
$ws.Cells.Item(70, 14).Value = -10540
$ws.Cells.Item(73, 8).Value = 7166.6665
$ws.Cells.Item(73, 9).Value = 5750
$ws.Cells.Item(73, 10).Value = 10000
$ws.Cells.Item(73, 11).Value = 5750
$ws.Cells.Item(73, 12).Value = 10000
$ws.Cells.Item(73, 13).Value = -4814
$ws.Cells.Item(73, 14).Value = -11872
$ws.Cells.Item(93, 8).Value = 0
$ws.Cells.Item(93, 10).Value = 0
$ws.Cells.Item(93, 12).Value = 0
$ws.Cells.Item(93, 14).ClearContents()
$ws.Cells.Item(126, 8).Value = 1699.75
$ws.Cells.Item(126, 9).Value = 1750
$ws.Cells.Item(126, 10).Value = 1649.5
$ws.Cells.Item(126, 11).Value = 5250
$ws.Cells.Item(126, 12).Value = 4948.5
$ws.Cells.Item(126, 13).Value = -2780
$ws.Cells.Item(126, 14).Value = -9888.5
$ws.Cells.Item(132, 8).Value = 4772.5
$ws.Cells.Item(132, 9).Value = 4545
$ws.Cells.Item(132, 10).Value = 5000
$ws.Cells.Item(132, 11).Value = 13635
$ws.Cells.Item(132, 12).Value = 15000
$ws.Cells.Item(132, 13).Value = -11105
$ws.Cells.Item(132, 14).Value = -20060

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(61, 8).Value = 3098.625
$ws.Cells.Item(61, 9).Value = 1548.1666
$ws.Cells.Item(61, 11).Value = 1548.1666
$ws.Cells.Item(61, 13).Value = -1346.1666
$ws.Cells.Item(68, 8).Value = 7593.143
$ws.Cells.Item(68, 9).Value = 4550.6665
$ws.Cells.Item(68, 10).Value = 9875
$ws.Cells.Item(68, 11).Value = 4550.6665
$ws.Cells.Item(68, 12).Value = 9875
$ws.Cells.Item(68, 13).Value = -3801.6665
$ws.Cells.Item(68, 14).Value = -11373
$ws.Cells.Item(71, 8).Value = 7593.143
$ws.Cells.Item(71, 9).Value = 4550.6665
$ws.Cells.Item(71, 10).Value = 9875
$ws.Cells.Item(71, 11).Value = 22753.3325
$ws.Cells.Item(71, 12).Value = 49375
$ws.Cells.Item(71, 13).Value = -19009.3325
$ws.Cells.Item(71, 14).Value = -56863
$ws.Cells.Item(113, 8).Value = 3098.625
$ws.Cells.Item(113, 9).Value = 1548.1666
$ws.Cells.Item(113, 11).Value = 1548.1666
$ws.Cells.Item(113, 13).Value = 621.8334

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(62, 8).Value = 8060.3076
$ws.Cells.Item(62, 9).Value = 6908.6665
$ws.Cells.Item(62, 10).Value = 9047.429
$ws.Cells.Item(62, 11).Value = 6908.6665
$ws.Cells.Item(62, 12).Value = 9047.429
$ws.Cells.Item(62, 13).Value = -6284.6665
$ws.Cells.Item(62, 14).Value = -10295.429
$ws.Cells.Item(65, 8).Value = 8060.3076
$ws.Cells.Item(65, 9).Value = 6908.6665
$ws.Cells.Item(65, 10).Value = 9047.429
$ws.Cells.Item(65, 11).Value = 34543.3325
$ws.Cells.Item(65, 12).Value = 45237.145
$ws.Cells.Item(65, 13).Value = -31423.3325
$ws.Cells.Item(65, 14).Value = -51477.145
$ws.Cells.Item(113, 8).Value = 779.9167
$ws.Cells.Item(113, 9).Value = 686
$ws.Cells.Item(113, 11).Value = 2058
$ws.Cells.Item(113, 13).Value = 112
$ws.Cells.Item(141, 8).Value = 500000
$ws.Cells.Item(141, 10).Value = 0
$ws.Cells.Item(141, 12).Value = 0
$ws.Cells.Item(141, 14).ClearContents()
